# Apply the price / volume(1h) refresh scraped on 2023-05-10, plus the
# Algorand / TheSandbox row re-order that came with this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ('Price') cells are stored as literal text in the sheet (e.g. '27.594.10'
# or '0.3623'). Excel's Range.Value setter auto-converts numeric-looking strings to
# real numbers, which would silently drop formatting such as trailing zeros. Prefixing
# the value with a leading apostrophe forces Excel to keep/store it as text, matching
# the original workbook's cell typing.

$ws.Range("D2").Value = "27.611.33"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "1.840.91"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "'312.31"
$ws.Range("E5").Value = "  -1.11%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").Value = "'0.4259"
$ws.Range("E7").Value = "  +0.44%  "

$ws.Range("D8").Value = "'0.3612"
$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("D9").Value = "'0.07301"
$ws.Range("E9").Value = "  +0.66%  "

$ws.Range("D10").Value = "'0.8754"
$ws.Range("E10").Value = "  -1.81%  "

$ws.Range("D11").Value = "'20.59"
$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").Value = "1.873.10"
$ws.Range("E12").Value = "  +4.07%  "

$ws.Range("D13").Value = "'5.325"
$ws.Range("E13").Value = "  -0.57%  "

$ws.Range("D14").Value = "'6.495"
$ws.Range("E14").Value = "  -1.32%  "

$ws.Range("D15").Value = "'0.06953"
$ws.Range("E15").Value = "  +1.40%  "

$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("D17").Value = "'79.22"
$ws.Range("E17").Value = "  +0.75%  "

$ws.Range("D18").Value = "'0.000008930"
$ws.Range("E18").Value = "  +0.89%  "

$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.14%  "

$ws.Range("D20").Value = "'15.32"
$ws.Range("E20").Value = "  -0.73%  "

$ws.Range("D21").Value = "27.566.88"
$ws.Range("E21").Value = "  -0.19%  "

$ws.Range("D22").Value = "'4.970"
$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("D23").Value = "'10.30"
$ws.Range("E23").Value = "  -2.65%  "

$ws.Range("D24").Value = "2.084.75"
$ws.Range("E24").Value = "  +2.31%  "

$ws.Range("D25").Value = "'1.987"
$ws.Range("E25").Value = "  -2.05%  "

$ws.Range("D26").Value = "'155.56"
$ws.Range("E26").Value = "  +0.78%  "

$ws.Range("D27").Value = "'18.52"
$ws.Range("E27").Value = "  -0.23%  "

$ws.Range("D28").Value = "'119.00"
$ws.Range("E28").Value = "  -0.27%  "

$ws.Range("D29").Value = "'5.192"
$ws.Range("E29").Value = "  -0.69%  "

$ws.Range("D30").Value = "'1.874"
$ws.Range("E30").Value = "  +1.89%  "

$ws.Range("D31").Value = "'0.08878"
$ws.Range("E31").Value = "  -0.31%  "

$ws.Range("D32").Value = "'0.7582"
$ws.Range("E32").Value = "  -2.55%  "

$ws.Range("D33").Value = "'2.946"
$ws.Range("E33").Value = "  -0.58%  "

$ws.Range("D34").Value = "'4.492"
$ws.Range("E34").Value = "  -1.58%  "

$ws.Range("D35").Value = "'1.121"
$ws.Range("E35").Value = "  +1.75%  "

$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("D37").Value = "'0.05426"
$ws.Range("E37").Value = "  +0.56%  "

$ws.Range("D38").Value = "'1.104"
$ws.Range("E38").Value = "  +0.08%  "

$ws.Range("D39").Value = "'0.01926"
$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("D40").Value = "'2.823"
$ws.Range("E40").Value = "  -0.43%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5069"
$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1658"
$ws.Range("E42").Value = "  +0.51%  "

$ws.Range("D43").Value = "'6.524"
$ws.Range("E43").Value = "  -4.89%  "

$ws.Range("D44").Value = "'8.364"
$ws.Range("E44").Value = "  +1.14%  "

$ws.Range("D45").Value = "'0.06546"
$ws.Range("E45").Value = "  -0.99%  "

$ws.Range("D46").Value = "'10.36"
$ws.Range("E46").Value = "  +0.12%  "

$ws.Range("D47").Value = "'105.82"
$ws.Range("E47").Value = "  +0.86%  "

$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("D49").Value = "'0.4620"
$ws.Range("E49").Value = "  -1.84%  "

$ws.Range("D50").Value = "'1.635"
$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("D51").Value = "'64.33"
$ws.Range("E51").Value = "  -0.16%  "
